# EHPpUC workbook bug fix
# - Electrolyzer Hydrogen Production per Unit Capacity
# Updates the Hydrogen Energy Density constant (lb -> correct BTU/lb value),
# adds explanatory notes from Meghan about the FCEV-adoption study behind the
# electrolyzer-capacity assumption, and leaves behind the formatting
# footprint of the fix (a spacer row + a block of pasted-in formatting) the
# same way the original author's Excel session produced it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Calculations sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Calculations")
$ws.Activate()

# Bug fix: correct BTU/lb value for hydrogen energy density
$ws.Range("A8").Value = 60920

# Insert a blank spacer row right after row 8 (pushes the "lb per metric
# ton" block and everything below it down by one row, and Excel's COM
# layer auto-adjusts every dependent formula/reference for us).
$ws.Rows("9:9").Insert()
$ws.Rows("8:9").RowHeight = 15.4

# Give the new spacer row some (invisible / no-fill) formatting, matching
# the handful of cells the author's paste left behind.
$ws.Range("A9:G9").Interior.Pattern = -4142
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("B9:G9").Font.Name = "Calibri"
$ws.Range("C9").Font.Underline = $true
$ws.Range("A9:G9").ClearContents()

# New explanatory notes tied to the citation, placed next to the
# lb-per-metric-ton / BTU-H2-per-MW block
$ws.Range("C12").Value = "The study they are citing uses 3 scenarios of Fuel Cell Electric Vehicle adoption. Then it calculates the amount of hydrogen needed to suppor those vehicles. Then it calculates the electrolyzer capacity needed to supply that hydrogen. "
$ws.Range("C13").Value = "So, I think it's fair, using the EPS assumptions of 24/7/365 operation, that smallest electrolyzer you would need to produce 1.39e10 annual Btu would be 1 MW."
$ws.Range("C14").Value = "No reason to think this would be different for Texas."

# Leftover (no-fill) formatting block to the right / below the table,
# matching the footprint left by the author's copy/paste in Excel.
$ws.Range("F1:J7").Interior.Pattern = -4142
$ws.Range("K2:O3").Interior.Pattern = -4142
$ws.Range("D15:I28").Interior.Pattern = -4142

$ws.Range("J6").Select()

# ---------------------------------------------------------------------
# About sheet - restore as the active sheet / selection
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("B6").Select()

# ---------------------------------------------------------------------
# EHPpUC sheet - remember last selection
# ---------------------------------------------------------------------
$wsEHP = $wb.Worksheets.Item("EHPpUC")
$wsEHP.Activate()
$wsEHP.Range("B2").Select()

# Restore About as the active/visible sheet
$wsAbout.Activate()
